# "form submission and localstorage"
#
# 1) The footer "date" placeholder (a datetimeFigureOut field reading
#    "06-11-2023") is refreshed to "09-11-2023" everywhere it is defined:
#    the slide master and all eleven slide layouts.
# 2) Slide 3 ("GAMEPLAY") gets its scoring rules rebalanced:
#      "You reach a score of 15,30,45,60….." -> "You reach a score of 20,40,60….."
#      "lives+=1"                             -> "score+=10"
#      "When score==25,50,75,100……"           -> "When score==40,80,120….……"
#      "-invincibility duration for 5-10s if i.e. no lives-=1 if you miss
#       in that duration"                     -> "-life+=1"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the date field text on the slide master + every layout.
# ---------------------------------------------------------------------
$oldDate = "06-11-2023"
$newDate = "09-11-2023"

$master = $p.SlideMaster
foreach ($shp in $master.Shapes) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    foreach ($shp in $layout.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Update the gameplay scoring bullets on slide 3.
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$contentShape = $slide3.Shapes.Item(4)   # "Content Placeholder 4"
$tr = $contentShape.TextFrame.TextRange

# Paragraph 6: "You reach a score of 15,30,45,60…..-lives+=1"
$para6 = $tr.Paragraphs(6, 1)
$para6.Runs(1).Text = "You reach a score of 20,40,60….."
$para6.Runs(3).Text = "score+=10"

# Paragraph 8: "When score==25,50,75,100……-invincibility duration ..."
$para8 = $tr.Paragraphs(8, 1)
$para8.Runs(1).Text = "When score==40,80,120….……"
$para8.Runs(2).Text = "-life+=1"
